# Update "Neg_Change" (sheet1) and "Pos_Change" (sheet2) tables with refreshed
# market data. Both sheets keep the same header row (row 1); data rows are
# replaced starting at row 2, and any leftover rows from the previous
# (longer/shorter) table are cleared.

$wb = $excel.ActiveWorkbook

$sheet1Data = @(
    ,@('ADANIENSOL', 931.85, 932.95, 915, 928.25, 974648, 2119876, -0.5402334853548038, 'ADANIENSOL')
    ,@('ICICIGI', 1891.8, 1901.5, 1855.5, 1859, 312123, 636816, -0.5098694128288234, 'ICICIGI')
    ,@('LICHSGFIN', 561.6, 572.35, 561.05, 567.5, 987991, 2247538, -0.5604118818013311, 'LICHSGFIN')
    ,@('GMRAIRPORT', 90.71, 91.16, 89.85, 90.13, 4874972, 12062943, -0.5958720852780287, 'GMRAIRPORT')
    ,@('HINDPETRO', 456.5, 461.35, 450.85, 452.9, 1709268, 3613774, -0.5270130340192829, 'HINDPETRO')
    ,@('IGL', 218.35, 218.58, 215.21, 216, 1204131, 2695300, -0.5532478759321783, 'IGL')
    ,@('NATIONALUM', 230, 230.8, 223.4, 224.69, 9669882, 21222011, -0.5443465748839731, 'NATIONALUM')
    ,@('AMBER', 8250, 8382, 8206, 8278.5, 151471, 311542, -0.5138023123687978, 'AMBER')
    ,@('GRANULES', 569, 574.5, 560.25, 568, 1270829, 2945411, -0.5685393311833221, 'GRANULES')
    ,@('KAYNES', 7110, 7183, 7053.5, 7098, 344930, 691656, -0.5012983332755011, 'KAYNES')
    ,@('TITAGARH', 905, 912.3, 899, 901, 464518, 1111859, -0.5822150110760447, 'TITAGARH')
)

$sheet2Data = @(
    ,@('MARUTI', 15980, 16315, 15911, 16285, 423031, 300499, 0.4077617562787231, 'MARUTI')
    ,@('WIPRO', 248, 251.25, 246.25, 248.98, 7070200, 4629760, 0.5271201962952723, 'WIPRO')
    ,@('BEL', 409, 415.5, 408.35, 412.15, 13813014, 9035150, 0.5288084868541197, 'BEL')
    ,@('GRASIM', 2810.6, 2832.8, 2788.3, 2812, 582587, 404859, 0.438987400551797, 'GRASIM')
    ,@('LODHA', 1130.3, 1161.3, 1120.6, 1153.6, 1061874, 682554, 0.5557362494396048, 'LODHA')
    ,@('HAVELLS', 1481.2, 1520, 1480.2, 1499.9, 817076, 561521, 0.4551120973213825, 'HAVELLS')
    ,@('CANBK', 126.5, 128.6, 126.35, 127.36, 23836453, 16981691, 0.4036560316637489, 'CANBK')
    ,@('TATAPOWER', 396, 396.8, 389.6, 390.5, 7725390, 5139999, 0.5029944558354973, 'TATAPOWER')
    ,@('INDUSINDBK', 749.6, 766.75, 747.75, 763.5, 5386722, 3506978, 0.5360010812728223, 'INDUSINDBK')
    ,@('PIIND', 3528, 3606.5, 3524, 3563.9, 335525, 211896, 0.5834418771472798, 'PIIND')
    ,@('GLENMARK', 1935.9, 1978.2, 1925.1, 1938.4, 867422, 596524, 0.4541275791083008, 'GLENMARK')
    ,@('ALKEM', 5480.5, 5500, 5412.5, 5460, 171470, 108874, 0.5749398387126403, 'ALKEM')
    ,@('PAYTM', 1252, 1257.6, 1230.5, 1237, 4322055, 3067327, 0.4090623529868188, 'PAYTM')
    ,@('MFSL', 1601, 1602.9, 1579.2, 1584, 528833, 342605, 0.543564746574043, 'MFSL')
    ,@('MUTHOOTFIN', 3233, 3240, 3150.5, 3165, 542445, 366182, 0.4813535345811645, 'MUTHOOTFIN')
    ,@('HFCL', 74.35, 77.67, 74.35, 76.01, 23042103, 14941670, 0.5421370569688663, 'HFCL')
    ,@('ANGELONE', 2264.5, 2319.4, 2253, 2303.8, 1026736, 723562, 0.4190021034824936, 'ANGELONE')
    ,@('NUVAMA', 7143.5, 7314.5, 7117.5, 7243, 144871, 90905, 0.5936527143721467, 'NUVAMA')
    ,@('LAURUSLABS', 873, 891.15, 843, 877.2, 3985430, 2609957, 0.5270098319627489, 'LAURUSLABS')
    ,@('ABFRL', 84.55, 85.18, 83, 83.2, 3558376, 2239381, 0.5889998173602438, 'ABFRL')
    ,@('PNBHOUSING', 877, 878, 851, 855, 1230300, 822722, 0.4954018489842255, 'PNBHOUSING')
)

function Update-Table {
    param($Worksheet, $Data)

    $rowCount = $Data.Count
    $colCount = $Data[0].Count

    # Clear out the previous data body (row 2 downwards) before writing the
    # fresh rows, so a shrinking table doesn't leave stale trailing rows.
    $usedRange = $Worksheet.UsedRange
    $lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
    if ($lastRow -ge 2) {
        $Worksheet.Range($Worksheet.Cells.Item(2, 1), $Worksheet.Cells.Item($lastRow, $colCount)).ClearContents()
    }

    # Write the new data body as a single block assignment starting at A2.
    $startRow = 2
    for ($r = 0; $r -lt $rowCount; $r++) {
        $rowData = $Data[$r]
        for ($c = 0; $c -lt $colCount; $c++) {
            $Worksheet.Cells.Item($startRow + $r, $c + 1).Value = $rowData[$c]
        }
    }
}

$wsNeg = $wb.Worksheets.Item("Neg_Change")
$wsPos = $wb.Worksheets.Item("Pos_Change")

Update-Table $wsNeg $sheet1Data
Update-Table $wsPos $sheet2Data
